$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking strings
# (e.g. "1.000", "0.9996", "45.42") are preserved exactly as text rather
# than being auto-converted to numbers by Excel.
$priceVolumeRange = $ws.Range("D2:E51")
$priceVolumeRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.580.13"
$ws.Range("E2").Value = "  -2.60%  "
$ws.Range("D3").Value = "1.747.66"
$ws.Range("E3").Value = "  -3.42%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "322.55"
$ws.Range("E5").Value = "  -4.58%  "
$ws.Range("D6").Value = "0.9989"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.4261"
$ws.Range("E7").Value = "  -8.78%  "
$ws.Range("D8").Value = "0.3613"
$ws.Range("E8").Value = "  -5.45%  "
$ws.Range("D9").Value = "45.42"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "0.07458"
$ws.Range("E10").Value = "  -2.48%  "
$ws.Range("D11").Value = "1.116"
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("D12").Value = "0.9963"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "21.51"
$ws.Range("E13").Value = "  -4.56%  "
$ws.Range("D14").Value = "6.107"
$ws.Range("E14").Value = "  -3.82%  "
$ws.Range("D15").Value = "7.210"
$ws.Range("E15").Value = "  -3.37%  "
$ws.Range("D16").Value = "1.741.23"
$ws.Range("E16").Value = "  -3.61%  "
$ws.Range("D17").Value = "0.00001068"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").Value = "87.69"
$ws.Range("E18").Value = "  +7.01%  "
$ws.Range("D19").Value = "0.06246"
$ws.Range("E19").Value = "  -7.00%  "
$ws.Range("D20").Value = "0.9996"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "16.93"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("D22").Value = "6.125"
$ws.Range("E22").Value = "  -4.77%  "
$ws.Range("D23").Value = "0.5246"
$ws.Range("E23").Value = "  -5.62%  "
$ws.Range("D24").Value = "27.567.89"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").Value = "11.63"
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("D26").Value = "2.322"
$ws.Range("E26").Value = "  -3.74%  "
$ws.Range("D27").Value = "20.47"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.364"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "151.40"
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("D30").Value = "1.935.81"
$ws.Range("E30").Value = "  -3.88%  "
$ws.Range("D31").Value = "1.224"
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("D32").Value = "126.71"
$ws.Range("E32").Value = "  -5.16%  "
$ws.Range("D33").Value = "5.703"
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("D34").Value = "0.09140"
$ws.Range("E34").Value = "  -4.52%  "
$ws.Range("D35").Value = "3.678"
$ws.Range("E35").Value = "  -8.85%  "
$ws.Range("D36").Value = "12.72"
$ws.Range("E36").Value = "  +4.80%  "
$ws.Range("D37").Value = "0.02301"
$ws.Range("E37").Value = "  -2.61%  "
$ws.Range("D38").Value = "0.2141"
$ws.Range("E38").Value = "  -5.91%  "
$ws.Range("D39").Value = "5.087"
$ws.Range("E39").Value = "  -3.84%  "
$ws.Range("D40").Value = "0.06089"
$ws.Range("E40").Value = "  -4.91%  "
$ws.Range("D41").Value = "0.6431"
$ws.Range("E41").Value = "  -3.46%  "
$ws.Range("D42").Value = "1.196"
$ws.Range("E42").Value = "  -3.33%  "
$ws.Range("D43").Value = "1.415"
$ws.Range("E43").Value = "  -5.13%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "7.944"
$ws.Range("E44").Value = "  -4.59%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "0.9984"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").Value = "13.76"
$ws.Range("E46").Value = "  -3.73%  "
$ws.Range("D47").Value = "3.728"
$ws.Range("E47").Value = "  -3.43%  "
$ws.Range("D48").Value = "0.5901"
$ws.Range("E48").Value = "  -4.27%  "
$ws.Range("D49").Value = "125.60"
$ws.Range("E49").Value = "  -4.08%  "
$ws.Range("D50").Value = "1.963"
$ws.Range("E50").Value = "  -4.00%  "
$ws.Range("D51").Value = "0.06868"
$ws.Range("E51").Value = "  -4.04%  "

# Restore the original (default) cell style now that the text values are set,
# so the cells keep matching their original unstyled appearance.
$priceVolumeRange.Style = "Normal"
